# Atualizando o arquivo XLSX
# Applies odds/score updates to rows 2-7 per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.1
$ws.Range("I2").Value = 2.05
$ws.Range("L2").Value = 2.88
$ws.Range("AH2").Value = 41
$ws.Range("G3").Value = 1.95
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 2.75
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("Q3").Value = 1.9
$ws.Range("R3").Value = 1.95
$ws.Range("Y3").Value = 1.57
$ws.Range("Z3").Value = 2.25
$ws.Range("AD3").Value = 8
$ws.Range("AE3").Value = 9.5
$ws.Range("AF3").Value = 17
$ws.Range("AN3").Value = 9
$ws.Range("AO3").Value = 21
$ws.Range("AP3").Value = 15
$ws.Range("G4").Value = 1.5
$ws.Range("I4").Value = 7.5
$ws.Range("J4").Value = 2.1
$ws.Range("L4").Value = 7.5
$ws.Range("N4").Value = 7.5
$ws.Range("AN4").Value = 15
$ws.Range("AR4").Value = 67
$ws.Range("G5").Value = 2.55
$ws.Range("I5").Value = 2.8
$ws.Range("J5").Value = 3.25
$ws.Range("Y5").Value = 1.44
$ws.Range("Z5").Value = 2.63
$ws.Range("AA5").Value = 1.8
$ws.Range("AB5").Value = 1.91
$ws.Range("AC5").Value = 8
$ws.Range("AG5").Value = 21
$ws.Range("AI5").Value = 9
$ws.Range("AM5").Value = 251
$ws.Range("AN5").Value = 8.5
$ws.Range("G6").Value = 3.25
$ws.Range("H6").Value = 3.25
$ws.Range("I6").Value = 2.25
$ws.Range("L6").Value = 2.88
$ws.Range("N6").Value = 9.5
$ws.Range("O6").Value = 1.3
$ws.Range("P6").Value = 3.4
$ws.Range("T6").Value = 1.8
$ws.Range("AC6").Value = 10
$ws.Range("AD6").Value = 17
$ws.Range("AI6").Value = 9.5
$ws.Range("G7").Value = 2.8
$ws.Range("H7").Value = 3.1
$ws.Range("I7").Value = 2.6
$ws.Range("J7").Value = 3.4
$ws.Range("L7").Value = 3.25
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 10
$ws.Range("O7").Value = 1.3
$ws.Range("P7").Value = 3.4
$ws.Range("S7").Value = 2.05
$ws.Range("T7").Value = 1.8
$ws.Range("W7").Value = 3.5
$ws.Range("X7").Value = 1.29
$ws.Range("Y7").Value = 1.44
$ws.Range("Z7").Value = 2.63
$ws.Range("AA7").Value = 1.8
$ws.Range("AB7").Value = 1.91
$ws.Range("AC7").Value = 8.5
$ws.Range("AD7").Value = 13
$ws.Range("AE7").Value = 10
$ws.Range("AF7").Value = 26
$ws.Range("AG7").Value = 21
$ws.Range("AH7").Value = 29
$ws.Range("AI7").Value = 9.5
$ws.Range("AJ7").Value = 6
$ws.Range("AL7").Value = 51
$ws.Range("AM7").Value = 251
$ws.Range("AN7").Value = 8.5
$ws.Range("AO7").Value = 13
$ws.Range("AP7").Value = 10
$ws.Range("AQ7").Value = 26
$ws.Range("AR7").Value = 21
$ws.Range("AS7").Value = 29
